# Zweite Version mit ActionListenern und erster Logik.
# Updates the "Menü" sheet: richer "Zutaten" descriptions, refreshed
# picture-hyperlink URLs, a new "Capricciosa" row, and top-aligned /
# wrap-text styling across the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menü")

$xlTop    = -4160
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) New hyperlink target text for the existing rows (Salami, Prosciutto,
#    Funghi, Tonno) -- written in this order so the shared-string table
#    grows the same way the original authoring session did.
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "https://kungfu-pizza.ro/wp-content/uploads/2020/12/pizza-salami.jpg"
$ws.Range("E4").Value = "https://as2.ftcdn.net/v2/jpg/04/80/82/95/1000_F_480829568_QySaJtZXCxPULUZ3CFiyDIrnuxQ85Of9.jpg"
$ws.Range("E5").Value = "https://www.eatbetter.de/sites/eatbetter.de/files/styles/facebook/public/2023-04/pizza_funghi_8823.jpg?h=4521fff0&itok=XcDJ5gFv"
$ws.Range("E6").Value = "https://www.globus.de/media/globus/rezepte/globus/pizza_tonno_169.jpg"

# ---------------------------------------------------------------------
# 2) Expanded "Zutaten" descriptions for the existing rows.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "mit frischen Tomaten, Basillikum und Käse"
$ws.Range("C3").Value = "mit Salami und Käse"
$ws.Range("C4").Value = "mit Schinken und Käse"
$ws.Range("C5").Value = "mit frischen Pilzen"
$ws.Range("C6").Value = "mit Thunfisch und Zwiebeln"

# ---------------------------------------------------------------------
# 3) New row 7: Capricciosa.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Capricciosa"
$ws.Range("C7").Value = "mit Artischocken, Pilzen, Schinken und schwarzen Oliven"
$ws.Range("D7").Value = 7.99

# ---------------------------------------------------------------------
# 4) Hyperlinks: drop the stale ones on Salami/Prosciutto/Funghi/Tonno
#    (Funghi and Tonno end up with no link at all), re-point
#    Salami/Prosciutto at their new image, and link the new Capricciosa
#    picture. Re-scan the live collection after every single delete --
#    this engine's Hyperlinks collection renumbers on delete and stale
#    handles silently no-op otherwise.
# ---------------------------------------------------------------------
function Remove-HyperlinkAtRow($sheet, $targetRow) {
    $again = $true
    while ($again) {
        $again = $false
        foreach ($h in $sheet.Hyperlinks) {
            if ($h.Range.Row -eq $targetRow) {
                $h.Delete()
                $again = $true
                break
            }
        }
    }
}

Remove-HyperlinkAtRow $ws 3
Remove-HyperlinkAtRow $ws 4
Remove-HyperlinkAtRow $ws 5
Remove-HyperlinkAtRow $ws 6

$ws.Hyperlinks.Add($ws.Range("E3"), "https://kungfu-pizza.ro/wp-content/uploads/2020/12/pizza-salami.jpg")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://as2.ftcdn.net/v2/jpg/04/80/82/95/1000_F_480829568_QySaJtZXCxPULUZ3CFiyDIrnuxQ85Of9.jpg")
$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.italianstylecooking.net/wp-content/uploads/2022/01/Pizza-capricciosa.jpg")

# ---------------------------------------------------------------------
# 5) Formatting: every data row gets vertical="top"; the "Zutaten"
#    column wraps its (now longer) text and stays unlocked.
# ---------------------------------------------------------------------
$ws.Range("A2").VerticalAlignment = $xlTop

$ws.Range("B2").VerticalAlignment = $xlTop

$ws.Range("C2").WrapText = $true
$ws.Range("C2").VerticalAlignment = $xlTop
$ws.Range("C2").Locked = $false

$ws.Range("D2").VerticalAlignment = $xlTop

$ws.Range("E2").VerticalAlignment = $xlTop

$ws.Range("A3").HorizontalAlignment = $xlCenter
$ws.Range("A3").VerticalAlignment = $xlTop

$ws.Range("B3:B7").VerticalAlignment = $xlTop
$ws.Range("D3:D7").VerticalAlignment = $xlTop
$ws.Range("E3:E7").VerticalAlignment = $xlTop

$ws.Range("C3:C7").WrapText = $true
$ws.Range("C3:C7").VerticalAlignment = $xlTop
$ws.Range("C3:C7").Locked = $false

$ws.Range("A4:A7").HorizontalAlignment = $xlCenter
$ws.Range("A4:A7").VerticalAlignment = $xlTop

# ---------------------------------------------------------------------
# 6) Row heights for the wrapped, multi-line descriptions.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 43.2

# ---------------------------------------------------------------------
# 7) Selection matches the author's last edit (price of the new row).
# ---------------------------------------------------------------------
$ws.Range("D7").Select()
